$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "List Bullet") {
        $p.Range.InsertBefore("Design: ")
    }
}
